# Apply the "Updated symbol list" crypto-price refresh to Sheet1.
# Values that look numeric are written with a leading apostrophe so Excel
# keeps them stored as text (matching the original inlineStr cell type)
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'249.93"
$ws.Range("G2").Value = "'20"
$ws.Range("D3").Value = "'21.98"
$ws.Range("G3").Value = "'20"
$ws.Range("G4").Value = "'20"
$ws.Range("D5").Value = "'0.05644"
$ws.Range("G5").Value = "'20"
$ws.Range("D6").Value = "'6.481"
$ws.Range("G6").Value = "'20"
$ws.Range("D7").Value = "'0.8029"
$ws.Range("G7").Value = "'20"
$ws.Range("G8").Value = "'20"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1437"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("G9").Value = "'20"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07307"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("G10").Value = "'20"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03107"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G11").Value = "'20"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02914"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("G12").Value = "'20"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09265"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("G13").Value = "'20"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001666"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("G14").Value = "'20"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.212"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Value = "'20"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04733"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("G16").Value = "'20"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005812"
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("G17").Value = "'20"
$ws.Range("D18").Value = "'0.006399"
$ws.Range("G18").Value = "'20"
$ws.Range("D19").Value = "'0.005069"
$ws.Range("G19").Value = "'20"
$ws.Range("D20").Value = "'0.001052"
$ws.Range("G20").Value = "'20"
$ws.Range("G21").Value = "'20"
$ws.Range("D22").Value = "'3.976"
$ws.Range("G22").Value = "'20"
$ws.Range("G23").Value = "'20"
$ws.Range("D24").Value = "'2.114"
$ws.Range("G24").Value = "'20"
$ws.Range("G25").Value = "'20"
$ws.Range("G26").Value = "'20"
$ws.Range("D27").Value = "'0.0002901"
$ws.Range("G27").Value = "'20"
$ws.Range("G28").Value = "'20"
$ws.Range("G29").Value = "'20"
$ws.Range("G30").Value = "'20"
$ws.Range("G31").Value = "'20"
$ws.Range("G32").Value = "'20"
$ws.Range("G33").Value = "'20"
$ws.Range("G34").Value = "'20"
$ws.Range("G35").Value = "'20"
$ws.Range("G36").Value = "'20"
$ws.Range("G37").Value = "'20"
$ws.Range("G38").Value = "'20"
$ws.Range("G39").Value = "'20"
$ws.Range("D40").Value = "'0.04164"
$ws.Range("G40").Value = "'20"
$ws.Range("D41").Value = "'0.007082"
$ws.Range("G41").Value = "'20"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("G42").Value = "'20"
$ws.Range("D43").Value = "'0.1043"
$ws.Range("G43").Value = "'20"
$ws.Range("D44").Value = "'0.009369"
$ws.Range("G44").Value = "'20"
$ws.Range("G45").Value = "'20"
$ws.Range("G46").Value = "'20"
$ws.Range("D47").Value = "'0.6802"
$ws.Range("G47").Value = "'20"
$ws.Range("D48").Value = "'0.01550"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("G48").Value = "'20"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("G49").Value = "'20"
$ws.Range("G50").Value = "'20"
$ws.Range("G51").Value = "'20"
